$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1259
$ws1.Range("F3").Value = 1152
$ws1.Range("F4").Value = 883
$ws1.Range("F5").Value = 102
$ws1.Range("F6").Value = 63
$ws1.Range("F7").Value = 642
$ws1.Range("F8").Value = 93
$ws1.Range("F9").Value = 46
$ws1.Range("F11").Value = 2283
$ws1.Range("F12").Value = 1571
$ws1.Range("F13").Value = 1300
$ws1.Range("F15").Value = 227
$ws1.Range("F16").Value = 522
$ws1.Range("F17").Value = 736
$ws1.Range("F18").Value = 26
$ws1.Range("F19").Value = 273
$ws1.Range("F20").Value = 1084
$ws1.Range("F22").Value = 6
$ws1.Range("F24").Value = 4421
$ws1.Range("F25").Value = 207
$ws1.Range("F26").Value = 15
$ws1.Range("F29").Value = 192
$ws1.Range("F30").Value = 79
$ws1.Range("F32").Value = 636
$ws1.Range("F35").Value = 37
$ws1.Range("G35").Value = 65
$ws1.Range("F37").Value = 361
$ws1.Range("F38").Value = 945
$ws1.Range("F39").Value = 121
$ws1.Range("F41").Value = 122
$ws1.Range("F42").Value = 109

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 423
$ws2.Range("F6").Value = 2
$ws2.Range("F9").Value = 2

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1259
$ws4.Range("F5").Value = 1152
$ws4.Range("F6").Value = 883
$ws4.Range("F8").Value = 423
$ws4.Range("F9").Value = 102
$ws4.Range("F10").Value = 63
$ws4.Range("F11").Value = 642
$ws4.Range("F12").Value = 93
$ws4.Range("F13").Value = 46
$ws4.Range("F14").Value = 2
$ws4.Range("F17").Value = 2283
$ws4.Range("F18").Value = 1571
$ws4.Range("F19").Value = 1300
$ws4.Range("F21").Value = 227
$ws4.Range("F22").Value = 522
$ws4.Range("F24").Value = 736
$ws4.Range("F25").Value = 26
$ws4.Range("F26").Value = 273
$ws4.Range("F27").Value = 1084
$ws4.Range("F28").Value = 6
$ws4.Range("F29").Value = 4422
$ws4.Range("F30").Value = 207
$ws4.Range("F31").Value = 15
$ws4.Range("F34").Value = 192
$ws4.Range("F35").Value = 79
$ws4.Range("F37").Value = 636
$ws4.Range("F40").Value = 37
$ws4.Range("G40").Value = 65
$ws4.Range("F41").Value = 361
$ws4.Range("F42").Value = 945
$ws4.Range("F43").Value = 121
$ws4.Range("F45").Value = 122
$ws4.Range("F46").Value = 109
$ws4.Range("F47").Value = 2
